$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6468.125
$ws.Range("I18").Value = 6807.6665
$ws.Range("K18").Value = 6807.6665
$ws.Range("M18").Value = -6523.6665
$ws.Range("H64").Value = 7495
$ws.Range("J64").Value = 7495
$ws.Range("L64").Value = 7495
$ws.Range("N64").Value = -7991
$ws.Range("H67").Value = 7495
$ws.Range("J67").Value = 7495
$ws.Range("L67").Value = 7495
$ws.Range("N67").Value = -9211
$ws.Range("H138").Value = 5580.528
$ws.Range("I138").Value = 2256.0715
$ws.Range("J138").Value = 7696.091
$ws.Range("K138").Value = 6768.2145
$ws.Range("L138").Value = 23088.273
$ws.Range("M138").Value = -1628.2145
$ws.Range("N138").Value = -33368.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3438.125
$ws.Range("I2").Value = 1687
$ws.Range("K2").Value = 1687
$ws.Range("M2").Value = -1574
$ws.Range("H32").Value = 2348.8193
$ws.Range("I32").Value = 2030.4492
$ws.Range("K32").Value = 2030.4492
$ws.Range("M32").Value = -1743.4492
$ws.Range("H61").Value = 5607.952
$ws.Range("I61").Value = 2626.2856
$ws.Range("J61").Value = 11571.286
$ws.Range("K61").Value = 2626.2856
$ws.Range("L61").Value = 11571.286
$ws.Range("M61").Value = -2414.2856
$ws.Range("N61").Value = -11995.286
$ws.Range("H74").Value = 43058.273
$ws.Range("I74").Value = 57825.38
$ws.Range("J74").Value = 4126.8184
$ws.Range("K74").Value = 57825.38
$ws.Range("L74").Value = 4126.8184
$ws.Range("M74").Value = -56951.38
$ws.Range("N74").Value = -5874.8184
$ws.Range("H77").Value = 43058.273
$ws.Range("I77").Value = 57825.38
$ws.Range("J77").Value = 4126.8184
$ws.Range("K77").Value = 289126.9
$ws.Range("L77").Value = 20634.092
$ws.Range("M77").Value = -284758.9
$ws.Range("N77").Value = -29370.092
$ws.Range("H102").Value = 2386.7778
$ws.Range("I102").Value = 2578.7144
$ws.Range("K102").Value = 2578.7144
$ws.Range("M102").Value = -956.7143999999998
$ws.Range("H116").Value = 3438.125
$ws.Range("I116").Value = 1687
$ws.Range("K116").Value = 1687
$ws.Range("M116").Value = 607
$ws.Range("H122").Value = 33156.285
$ws.Range("I122").Value = 67533.336
$ws.Range("J122").Value = 7373.5
$ws.Range("K122").Value = 202600.008
$ws.Range("L122").Value = 22120.5
$ws.Range("M122").Value = -200150.008
$ws.Range("N122").Value = -27020.5
$ws.Range("H136").Value = 5607.952
$ws.Range("I136").Value = 2626.2856
$ws.Range("J136").Value = 11571.286
$ws.Range("K136").Value = 7878.8568
$ws.Range("L136").Value = 34713.858
$ws.Range("M136").Value = -5328.8568
$ws.Range("N136").Value = -39813.858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3438.125
$ws.Range("I3").Value = 1687
$ws.Range("K3").Value = 1687
$ws.Range("M3").Value = -1573
$ws.Range("H86").Value = 69503600
$ws.Range("I86").Value = 35859756
$ws.Range("J86").Value = 90913310
$ws.Range("K86").Value = 35859756
$ws.Range("L86").Value = 90913310
$ws.Range("M86").Value = -35858633
$ws.Range("N86").Value = -90915556
$ws.Range("H89").Value = 69503600
$ws.Range("I89").Value = 35859756
$ws.Range("J89").Value = 90913310
$ws.Range("K89").Value = 179298780
$ws.Range("L89").Value = 454566550
$ws.Range("M89").Value = -179293164
$ws.Range("N89").Value = -454577782
$ws.Range("H105").Value = 2616.7144
$ws.Range("I105").Value = 2037.4667
$ws.Range("J105").Value = 4064.8333
$ws.Range("K105").Value = 2037.4667
$ws.Range("L105").Value = 4064.8333
$ws.Range("M105").Value = -290.4666999999999
$ws.Range("N105").Value = -7558.8333
$ws.Range("H134").Value = 4324.2383
$ws.Range("I134").Value = 1646.8718
$ws.Range("J134").Value = 8674.958
$ws.Range("K134").Value = 4940.6154
$ws.Range("L134").Value = 26024.874
$ws.Range("M134").Value = -2405.6154
$ws.Range("N134").Value = -31094.874

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 15878500
$ws.Range("I86").Value = 15878500
$ws.Range("K86").Value = 15878500
$ws.Range("M86").Value = -15877377
$ws.Range("H89").Value = 15878500
$ws.Range("I89").Value = 15878500
$ws.Range("K89").Value = 79392500
$ws.Range("M89").Value = -79386884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3334271.5
$ws.Range("J12").Value = 4546316.5
$ws.Range("L12").Value = 13638949.5
$ws.Range("N12").Value = -13639295.5
$ws.Range("H117").Value = 800.3
$ws.Range("I117").Value = 613
$ws.Range("J117").Value = 847.125
$ws.Range("K117").Value = 1839
$ws.Range("L117").Value = 2541.375
$ws.Range("M117").Value = 1603
$ws.Range("N117").Value = -9425.375
$ws.Range("H129").Value = 23953644
$ws.Range("I129").Value = 1111
$ws.Range("K129").Value = 3333
$ws.Range("M129").Value = 1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5357.0435
$ws.Range("I113").Value = 2270.6086
$ws.Range("J113").Value = 8443.479
$ws.Range("K113").Value = 2270.6086
$ws.Range("L113").Value = 8443.479
$ws.Range("M113").Value = -100.6086
$ws.Range("N113").Value = -12783.479
$ws.Range("H126").Value = 4406
$ws.Range("I126").Value = 4406
$ws.Range("K126").Value = 13218
$ws.Range("M126").Value = -10748
$ws.Range("H132").Value = 13903.5625
$ws.Range("I132").Value = 5911.2
$ws.Range("K132").Value = 17733.6
$ws.Range("M132").Value = -15203.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 23668.666
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H46").Value = 7409473.5
$ws.Range("J46").Value = 7409473.5
$ws.Range("L46").Value = 7409473.5
$ws.Range("N46").Value = -7409849.5
$ws.Range("H59").Value = 56447.332
$ws.Range("J59").Value = 56447.332
$ws.Range("L59").Value = 56447.332
$ws.Range("N59").Value = -57755.332
$ws.Range("H69").Value = 45000
$ws.Range("J69").Value = 45000
$ws.Range("L69").Value = 45000
$ws.Range("N69").Value = -46622
$ws.Range("H72").Value = 45000
$ws.Range("J72").Value = 45000
$ws.Range("L72").Value = 135000
$ws.Range("N72").Value = -143112
$ws.Range("H74").Value = 44444
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 44444
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H101").Value = 53531.332
$ws.Range("J101").Value = 53531.332
$ws.Range("L101").Value = 53531.332
$ws.Range("N101").Value = -60021.332
$ws.Range("H122").Value = 3995.1785
$ws.Range("I122").Value = 3394.7827
$ws.Range("J122").Value = 6757
$ws.Range("K122").Value = 10184.3481
$ws.Range("L122").Value = 20271
$ws.Range("M122").Value = -7734.348100000001
$ws.Range("N122").Value = -25171
$ws.Range("H132").Value = 13166282
$ws.Range("I132").Value = 22732758
$ws.Range("J132").Value = 12375
$ws.Range("K132").Value = 68198274
$ws.Range("L132").Value = 37125
$ws.Range("M132").Value = -68195744
$ws.Range("N132").Value = -42185
$ws.Range("H136").Value = 11845.435
$ws.Range("I136").Value = 3485
$ws.Range("J136").Value = 13099.5
$ws.Range("K136").Value = 10455
$ws.Range("L136").Value = 39298.5
$ws.Range("M136").Value = -7905
$ws.Range("N136").Value = -44398.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 24998
$ws.Range("I15").Value = 24998
$ws.Range("K15").Value = 24998
$ws.Range("M15").Value = -24710
$ws.Range("H113").Value = 1660.8518
$ws.Range("I113").Value = 1153.4
$ws.Range("K113").Value = 3460.2
$ws.Range("M113").Value = -1290.2
$ws.Range("H114").Value = 42870.8
$ws.Range("J114").Value = 42870.8
$ws.Range("L114").Value = 42870.8
$ws.Range("N114").Value = -51548.8
$ws.Range("H132").Value = 13170106
$ws.Range("I132").Value = 20839160
$ws.Range("J132").Value = 23156.143
$ws.Range("K132").Value = 62517480
$ws.Range("L132").Value = 69468.429
$ws.Range("M132").Value = -62514950
$ws.Range("N132").Value = -74528.429
